# Applies targeted cell-value corrections across multiple sheets
# (Leve profit/price recalculations) as captured in the commit diff.
$wb = $excel.ActiveWorkbook

# ALC!row44
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H44").Value = 14200
$ws.Range("J44").Value = 14200
$ws.Range("L44").Value = 14200
$ws.Range("N44").Value = -15124

# ALC!row70
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 2713.7222
$ws.Range("I70").Value = 883.63635
$ws.Range("J70").Value = 5589.5713
$ws.Range("K70").Value = 2650.90905
$ws.Range("L70").Value = 16768.7139
$ws.Range("M70").Value = -2380.90905
$ws.Range("N70").Value = -17308.7139

# ALC!row73
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 2713.7222
$ws.Range("I73").Value = 883.63635
$ws.Range("J73").Value = 5589.5713
$ws.Range("K73").Value = 2650.90905
$ws.Range("L73").Value = 16768.7139
$ws.Range("M73").Value = -1714.90905
$ws.Range("N73").Value = -18640.7139

# ALC!row132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2058348.5
$ws.Range("I132").Value = 760.1731
$ws.Range("J132").Value = 55555644
$ws.Range("K132").Value = 2280.5193
$ws.Range("L132").Value = 166666932
$ws.Range("M132").Value = 249.4807000000001
$ws.Range("N132").Value = -166671992

# ARM!row61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3865871.5
$ws.Range("I61").Value = 1667565.2
$ws.Range("J61").Value = 58823530
$ws.Range("K61").Value = 1667565.2
$ws.Range("L61").Value = 58823530
$ws.Range("M61").Value = -1667353.2
$ws.Range("N61").Value = -58823954

# ARM!row136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3865871.5
$ws.Range("I136").Value = 1667565.2
$ws.Range("J136").Value = 58823530
$ws.Range("K136").Value = 5002695.6
$ws.Range("L136").Value = 176470590
$ws.Range("M136").Value = -5000145.6
$ws.Range("N136").Value = -176475690

# CRP!row16
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 757.8125
$ws.Range("I16").Value = 732
$ws.Range("J16").Value = 814.6
$ws.Range("K16").Value = 732
$ws.Range("L16").Value = 814.6
$ws.Range("M16").Value = -445
$ws.Range("N16").Value = -1388.6

# CRP!row31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2267598.2
$ws.Range("I31").Value = 1303569.1
$ws.Range("J31").Value = 4471093
$ws.Range("K31").Value = 1303569.1
$ws.Range("L31").Value = 4471093
$ws.Range("M31").Value = -1303274.1
$ws.Range("N31").Value = -4471683

# CRP!row34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2267598.2
$ws.Range("I34").Value = 1303569.1
$ws.Range("J34").Value = 4471093
$ws.Range("K34").Value = 1303569.1
$ws.Range("L34").Value = 4471093
$ws.Range("M34").Value = -1303367.1
$ws.Range("N34").Value = -4471497

# CRP!row57
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H57").Value = 41000
$ws.Range("I57").Value = 13000
$ws.Range("J57").Value = 55000
$ws.Range("K57").Value = 13000
$ws.Range("L57").Value = 55000
$ws.Range("M57").Value = -12440
$ws.Range("N57").Value = -56120

# CRP!row113
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 757.8125
$ws.Range("I113").Value = 732
$ws.Range("J113").Value = 814.6
$ws.Range("K113").Value = 732
$ws.Range("L113").Value = 814.6
$ws.Range("M113").Value = 1438
$ws.Range("N113").Value = -5154.6

# CRP!row132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3574851.2
$ws.Range("I132").Value = 5558045
$ws.Range("J132").Value = 5102.4
$ws.Range("K132").Value = 16674135
$ws.Range("L132").Value = 15307.2
$ws.Range("M132").Value = -16671605
$ws.Range("N132").Value = -20367.2

# CRP!row134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1030215.5
$ws.Range("I134").Value = 5372.304
$ws.Range("J134").Value = 2503427.5
$ws.Range("K134").Value = 16116.912
$ws.Range("L134").Value = 7510282.5
$ws.Range("M134").Value = -13581.912
$ws.Range("N134").Value = -7515352.5

# CUL!row46
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 2040.7222
$ws.Range("J46").Value = 2040.7222
$ws.Range("L46").Value = 6122.1666
$ws.Range("N46").Value = -6304.1666

# CUL!row86
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 335
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 335
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 1005
$ws.Range("N86").Value = -3377
$ws.Range("M86").ClearContents()

# CUL!row89
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H89").Value = 335
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 335
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 3015
$ws.Range("N89").Value = -14871
$ws.Range("M89").ClearContents()

# CUL!row92
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 6097861
$ws.Range("I92").Value = 300
$ws.Range("J92").Value = 8130381.5
$ws.Range("K92").Value = 900
$ws.Range("L92").Value = 24391144.5
$ws.Range("M92").Value = 348
$ws.Range("N92").Value = -24393640.5

# CUL!row122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 828.9643
$ws.Range("I122").Value = 258.26315
$ws.Range("J122").Value = 2033.7778
$ws.Range("K122").Value = 2324.36835
$ws.Range("L122").Value = 18304.0002
$ws.Range("M122").Value = 125.6316500000003
$ws.Range("N122").Value = -23204.0002

# GSM!row58
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 8196
$ws.Range("I58").Value = 980
$ws.Range("K58").Value = 980
$ws.Range("M58").Value = -703

# GSM!row102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 5303.6816
$ws.Range("I102").Value = 5820.647
$ws.Range("J102").Value = 3546
$ws.Range("K102").Value = 5820.647
$ws.Range("L102").Value = 3546
$ws.Range("M102").Value = -4198.647
$ws.Range("N102").Value = -6790

# LTW!row132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2302163.2
$ws.Range("I132").Value = 3031692.2
$ws.Range("J132").Value = 9358.143
$ws.Range("K132").Value = 9095076.600000001
$ws.Range("L132").Value = 28074.429
$ws.Range("M132").Value = -9092546.600000001
$ws.Range("N132").Value = -33134.429

# LTW!row136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1177400
$ws.Range("I136").Value = 1681520.9
$ws.Range("J136").Value = 1117.9333
$ws.Range("K136").Value = 5044562.699999999
$ws.Range("L136").Value = 3353.7999
$ws.Range("M136").Value = -5042012.699999999
$ws.Range("N136").Value = -8453.7999

# WVR!row39
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 2933.3333
$ws.Range("I39").Value = 2933.3333
$ws.Range("K39").Value = 2933.3333
$ws.Range("M39").Value = -2520.3333

# WVR!row47
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()

# WVR!row52
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 5000
$ws.Range("J52").Value = 5000
$ws.Range("L52").Value = 5000
$ws.Range("M52").Value = -5452

# WVR!row126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 62507900
$ws.Range("I126").Value = 125000300
$ws.Range("J126").Value = 15499.5
$ws.Range("K126").Value = 375000900
$ws.Range("L126").Value = 46498.5
$ws.Range("M126").Value = -374998430
$ws.Range("N126").Value = -51438.5

# WVR!row132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1965050.5
$ws.Range("I132").Value = 1327984.9
$ws.Range("K132").Value = 3983954.7
$ws.Range("M132").Value = -3981424.7

Write-Output "Applied all Leve Profits cell updates."
